$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2..62 (columns A: index, B: date serial, C: value)
$data = @(
    @(2, 0, 44347, 1.181499232696965),
    @(3, 1, 44348, 1.146181993214942),
    @(4, 2, 44349, 1.150018702939569),
    @(5, 3, 44350, 0.8737762439103937),
    @(6, 4, 44351, 1.037960406863995),
    @(7, 5, 44352, 0.985850200243963),
    @(8, 6, 44353, 1.458273020860609),
    @(9, 7, 44354, 1.137700295024511),
    @(10, 8, 44355, 1.135817194695616),
    @(11, 9, 44356, 1.153391533819905),
    @(12, 10, 44357, 1.305306929603892),
    @(13, 11, 44358, 1.333599803680151),
    @(14, 12, 44359, 0.9793663816800745),
    @(15, 13, 44360, 1.365812944567083),
    @(16, 14, 44361, 1.230596766086818),
    @(17, 15, 44362, 1.290366449743396),
    @(18, 16, 44363, 1.327286889486378),
    @(19, 17, 44364, 1.400332457112883),
    @(20, 18, 44365, 1.402286566111398),
    @(21, 19, 44366, 1.128240833611947),
    @(22, 20, 44367, 0.787236574999207),
    @(23, 21, 44368, 1.316907245010859),
    @(24, 22, 44369, 1.316482349342449),
    @(25, 23, 44370, 1.20861428969807),
    @(26, 24, 44371, 1.186603379411141),
    @(27, 25, 44372, 1.272861131983207),
    @(28, 26, 44373, 1.202806977987106),
    @(29, 27, 44374, 0.8447195786317578),
    @(30, 28, 44375, 1.241982572455842),
    @(31, 29, 44376, 1.337140146150796),
    @(32, 30, 44377, 1.119982378331922),
    @(33, 31, 44378, 1.324395868968712),
    @(34, 32, 44379, 1.240457915831663),
    @(35, 33, 44380, 1.308216103096456),
    @(36, 34, 44381, 0.8398474682411881),
    @(37, 35, 44382, 0.9383992795788307),
    @(38, 36, 44383, 0.9623773872829288),
    @(39, 37, 44384, 0.9064416255117055),
    @(40, 38, 44385, 1.059515742713512),
    @(41, 39, 44386, 1.112029127187046),
    @(42, 40, 44387, 1.073792550223214),
    @(43, 41, 44388, 1.278090059703807),
    @(44, 42, 44389, 0.933675046958532),
    @(45, 43, 44390, 0.9951211744308759),
    @(46, 44, 44391, 0.9412997161814904),
    @(47, 45, 44392, 0.9055637335077964),
    @(48, 46, 44393, 1.148156867911272),
    @(49, 47, 44394, 1.187548089091369),
    @(50, 48, 44395, 1.123674076106842),
    @(51, 49, 44396, 0.6798248109855467),
    @(52, 50, 44397, 0.7185852435530086),
    @(53, 51, 44398, 0.8087669920013262),
    @(54, 52, 44399, 0.7263670793501542),
    @(55, 53, 44400, 0.950490201674277),
    @(56, 54, 44401, 0.9901610036375662),
    @(57, 55, 44402, 0.740482944011865),
    @(58, 56, 44403, 0.7003358898006764),
    @(59, 57, 44404, 0.7159132830279076),
    @(60, 58, 44405, 0.7643231395625296),
    @(61, 59, 44406, 0.8641530034235917),
    @(62, 60, 44407, 1.045156264862147)
)

$lastExistingRow = 30

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]

    if ($r -gt $lastExistingRow) {
        # New rows need the same formatting as the existing data rows:
        # Column A: bold, centered, thin border (matches style used in A2:A30)
        $aCell = $ws.Cells.Item($r, 1)
        $aCell.Font.Bold = $true
        $aCell.HorizontalAlignment = -4108
        $aCell.VerticalAlignment = -4160
        $aCell.Borders.LineStyle = 1
        $aCell.Borders.Weight = 2

        # Column B: custom date/time number format (matches style used in B2:B30)
        $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
}

Write-Output "updated rows 2-62"
